# Reposition the Meteorological influence/conditions tables and fix the
# alignment of the site-location sub-line on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# 1) "Testing Organization and Site Information" table: left-align the
#    "(City, State, Latitude and Longitude)" sub-line instead of centering it.
$siteInfoTable = Get-ShapeById $s 18
$tbl = $siteInfoTable.Table

$cell = $tbl.Rows.Item(3).Cells.Item(1)
$tr = $cell.Shape.TextFrame.TextRange
$secondLine = $tr.Paragraphs(2, 1)
if ($secondLine.Text -like "*City, State, Latitude*") {
    $secondLine.ParagraphFormat.Alignment = 1
}

# 2) Move the two small "Meteorological influence/conditions" tables.
$tempTable = Get-ShapeById $s 32
$tempTable.Left = 126.0
$tempTable.Top = 1439.28

$concTable = Get-ShapeById $s 33
$concTable.Left = 622.4697637795275
$concTable.Top = 1439.28
